# Rewrites the four bullets under "KEY ACHIEVEMENTS AND IMPACT" / "Impact" into
# six bullets, adding the Supreme Court mentions (per commit message) and swapping
# in the new "algorithmic innovation" / "breakthrough demographic discovery" framing.
#
# Before (paragraphs, in order):
#   • Discovered systematic race coding errors affecting all Black and Asian-American voters
#   • Algorithm reduced mapping costs by **73.5%**, saving campaigns and organizations **$4.7M**
#   • Built redistricting platform used by thousands of analysts nationwide
#   • Achieved **87%** prediction accuracy for voter turnout vs. industry standard of **71%**
#
# After (paragraphs, in order):
#   • Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs **73.5%**
#   • **$4.7M** savings enabled nonprofit access
#   • Legal precedent: Data analysis utilized in Supreme Court case
#   • Expert methodology validated at highest judicial level
#   • Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions
#   • **178%** accuracy improvement in racial classification algorithms

$d = $word.ActiveDocument

# Locate the first of the four existing bullets unambiguously: walk paragraphs
# looking for the "KEY ACHIEVEMENTS AND IMPACT" heading, then its "Impact" sub-
# heading, then the bullet right after it. (Plain substring Find() is unsafe here
# because near-identical bullet text also appears earlier, under "PROFESSIONAL
# EXPERIENCE".)
$firstBulletIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd() -eq "KEY ACHIEVEMENTS AND IMPACT") {
        $firstBulletIndex = $i + 2
        break
    }
}
if ($firstBulletIndex -eq 0) {
    throw "Could not locate 'KEY ACHIEVEMENTS AND IMPACT' heading"
}

# Sanity-check the four paragraphs we expect to be there.
$p1 = $d.Paragraphs.Item($firstBulletIndex).Range.Text
$p2 = $d.Paragraphs.Item($firstBulletIndex + 1).Range.Text
$p3 = $d.Paragraphs.Item($firstBulletIndex + 2).Range.Text
$p4 = $d.Paragraphs.Item($firstBulletIndex + 3).Range.Text
if ($p1 -notlike "*Discovered systematic race coding errors*") {
    throw "Unexpected paragraph 1: $p1"
}
if ($p2 -notlike "*Algorithm reduced mapping costs*") {
    throw "Unexpected paragraph 2: $p2"
}
if ($p3 -notlike "*Built redistricting platform*") {
    throw "Unexpected paragraph 3: $p3"
}
if ($p4 -notlike "*Achieved*") {
    throw "Unexpected paragraph 4: $p4"
}

# The 3rd paragraph ("Built redistricting platform...") is a single plain run with
# no rPr at all -- use it as a clean template. Duplicate it five more times right
# after itself, giving six consecutive "clean" paragraphs to fill in (this avoids
# ambiguous-formatting-at-insertion-point problems that arise from inserting plain
# text directly adjacent to bold runs).
$templateIndex = $firstBulletIndex + 2
$templateFormatted = $d.Paragraphs.Item($templateIndex).Range.FormattedText
for ($i = 0; $i -lt 5; $i++) {
    $insertPos = $d.Paragraphs.Item($templateIndex + $i).Range.End
    $insertRange = $d.Range($insertPos, $insertPos)
    $insertRange.FormattedText = $templateFormatted
}

# Fill in the six clean paragraphs with the target text, using placeholder tokens
# where a bold/colored run belongs.
$targets = @(
    "• Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs **BOLDTOKEN1**",
    "• **BOLDTOKEN2** savings enabled nonprofit access",
    "• Legal precedent: Data analysis utilized in Supreme Court case",
    "• Expert methodology validated at highest judicial level",
    "• Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions",
    "• **BOLDTOKEN3** accuracy improvement in racial classification algorithms"
)
for ($i = 0; $i -lt 6; $i++) {
    $d.Paragraphs.Item($templateIndex + $i).Range.Text = $targets[$i]
}

# Replace each placeholder token with its bold/colored value (matches the existing
# bold run styling used throughout this resume: Bold + RGB 2C3E50).
function Set-BoldToken($paraIndex, $token, $value) {
    $paraRange = $d.Paragraphs.Item($paraIndex).Range
    $find = $paraRange.Find
    $find.ClearFormatting()
    $find.Text = $token
    $find.Replacement.ClearFormatting()
    $find.Replacement.Font.Bold = 1
    $find.Replacement.Font.Color = 5258796   # RGB(0x2C,0x3E,0x50) packed as 0x00503E2C
    $find.Replacement.Text = $value
    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 0, $true, $find.Replacement.Text, 2) | Out-Null
}

Set-BoldToken $templateIndex "BOLDTOKEN1" "73.5%"
Set-BoldToken ($templateIndex + 1) "BOLDTOKEN2" "`$4.7M"
Set-BoldToken ($templateIndex + 5) "BOLDTOKEN3" "178%"

# Finally, delete the three original paragraphs that are no longer needed (delete
# highest index first so the lower indices stay valid):
#   - old 4th bullet ("Achieved 87% ... 71%"), now sitting right after our six new ones
#   - old 2nd bullet ("Algorithm reduced mapping costs by ...")
#   - old 1st bullet ("Discovered systematic race coding errors ...")
$d.Paragraphs.Item($templateIndex + 6).Range.Delete()
$d.Paragraphs.Item($firstBulletIndex + 1).Range.Delete()
$d.Paragraphs.Item($firstBulletIndex).Range.Delete()
